# "nambah deskripsi dan barang jenis" - add a new "deskripsi" header column
# (H1) with the same centered header styling as the rest of the header row,
# widen that column, drop a lone-space value down at I7, and leave the
# active selection on F5 (matching the author's saved state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell: H1 = "deskripsi", centered like the other header cells.
$ws.Range("H1").Value = "deskripsi"
$ws.Range("H1").HorizontalAlignment = -4108  # xlCenter

# Widen the new column to fit the longer header text.
$ws.Columns.Item(8).ColumnWidth = 19.83

# Stray single-space entry further down the sheet.
$ws.Range("I7").Value = " "

# Leave the active selection where the author left it.
$ws.Range("F5").Select() | Out-Null
